# stagingLog.xlsx - "updates to global missing page test"
#
# Row 3, column A previously held the literal text "test". The update
# replaces it with the same build string already recorded in A2 (e.g.
# "27c37e5e4b built at 2020-08-20 16:02"), so both data rows report the
# build that was used for the run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value2 = $ws.Range("A2").Value2

# Re-fit the row after the multi-line build string lands in A3 so the row
# keeps its default (non-custom) height, matching the original layout.
$ws.Rows(3).AutoFit()
